# moving calculate_euclidean_dist to child class and build constructor from excel to child, car and school
# Update the "address" column (D) on the "Child" sheet with new coordinate pairs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Child")

$ws.Range("D2").Value  = "-7,-5"
$ws.Range("D3").Value  = "1,-7"
$ws.Range("D4").Value  = "-7,6"
$ws.Range("D5").Value  = "-7,-8"
$ws.Range("D6").Value  = "0,-4"
$ws.Range("D7").Value  = "-6,3"
$ws.Range("D8").Value  = "-10,-2"
$ws.Range("D9").Value  = "-2,8"
$ws.Range("D10").Value = "9,-4"
$ws.Range("D11").Value = "-1,-8"
$ws.Range("D12").Value = "-5,-9"
$ws.Range("D13").Value = "-2,4"
$ws.Range("D14").Value = "2,7"
$ws.Range("D15").Value = "-7,6"
$ws.Range("D16").Value = "-9,-4"
$ws.Range("D17").Value = "7,6"
$ws.Range("D18").Value = "0,4"
$ws.Range("D19").Value = "-8,-2"
$ws.Range("D20").Value = "6,1"
$ws.Range("D21").Value = "1,8"
